# Apply "hybrid bold + color highlighting" to quantitative impact metrics
# (percentages, +/- margins, dollar amounts) inside selected bullet /
# paragraph runs, matching the target diff exactly: each highlighted
# number becomes its own run with <w:b/> and <w:color w:val="2C3E50"/>.

$d = $word.ActiveDocument

# wdColor values are packed 0x00BBGGRR, so build the constant from RGB.
function RGBColor([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}
$accentColor = RGBColor 44 62 80   # hex 2C3E50

# Locate the (first) paragraph whose text contains $needle.
function Get-ParaByText([string]$needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# Within a paragraph, find literal $text and apply bold + accent color
# to just that run of characters (Word automatically splits the
# surrounding run(s) so formatting only covers the matched text).
function Set-MetricHighlight($para, [string]$text) {
    $r = $para.Range.Duplicate
    $found = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $r.Font.Bold = 1
        $r.Font.Color = $accentColor
    }
    return $found
}

$pm = [char]0x00B1   # '±'
$pm42 = $pm + "4.2%"
$pm21 = $pm + "2.1%"
$dollar47M = "`$4.7M"
$dollar2 = "`$2"

# --- "Data Science & Political Analytics" bullets (Siege Analytics) -------

$p = Get-ParaByText "Discovered systematic race coding errors"
Set-MetricHighlight $p "23%" | Out-Null
Set-MetricHighlight $p "64%" | Out-Null

$p = Get-ParaByText "Utilized advanced sampling methods to decrease survey margin of error"
Set-MetricHighlight $p $pm42 | Out-Null
Set-MetricHighlight $p $pm21 | Out-Null
Set-MetricHighlight $p "71%" | Out-Null
Set-MetricHighlight $p "87%" | Out-Null

$p = Get-ParaByText "Trigonometric algorithm for boundary estimation"
Set-MetricHighlight $p "73.5%" | Out-Null
Set-MetricHighlight $p $dollar47M | Out-Null

$p = Get-ParaByText "Built real-time FEC analysis systems"
Set-MetricHighlight $p $dollar2 | Out-Null

# --- "Data Products Manager" bullet (Helm/Murmuration) --------------------

$p = Get-ParaByText "Modernized legacy ETL processes"
Set-MetricHighlight $p "57%" | Out-Null

# --- "KEY ACHIEVEMENTS AND IMPACT" bullets ---------------------------------

$p = Get-ParaByText "Predictive excellence: Utilized advanced sampling methods"
Set-MetricHighlight $p $pm42 | Out-Null
Set-MetricHighlight $p $pm21 | Out-Null

$p = Get-ParaByText "Increased voter turnout prediction accuracy"
Set-MetricHighlight $p "71%" | Out-Null
Set-MetricHighlight $p "87%" | Out-Null

$p = Get-ParaByText "Methodological advancement: Improved segmentation accuracy"
Set-MetricHighlight $p "34%" | Out-Null
Set-MetricHighlight $p "28%" | Out-Null
